$d = $word.ActiveDocument

$replacements = @(
    @("2025-12-22 Monday", "2025-12-23 Tuesday"),
    @("333×3=", "391×5="),
    @("708×8=", "329×3="),
    @("999×3=", "767×3="),
    @("119×8=", "267×2="),
    @("641×6=", "281×7="),
    @("938×9=", "144×8="),
    @("367×2=", "463×8="),
    @("909×4=", "321×9="),
    @("110×2=", "595×3="),
    @("921×5=", "363×2="),
    @("648×9=", "629×6="),
    @("609×6=", "139×7="),
    @("598×4=", "682×3="),
    @("499×7=", "719×2="),
    @("261×3=", "263×7="),
    @("819×6=", "301×4="),
    @("282×9=", "139×6="),
    @("488×6=", "383×3="),
    @("589×4=", "734×2="),
    @("960×2=", "733×7="),
    @("911×9=", "669×4="),
    @("907×4=", "635×9="),
    @("642×5=", "423×5="),
    @("666×5=", "387×4="),
    @("737×3=", "849×3=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
